# Refresh scraped "想去人数" (interest-count) figures in column F
# across the 展览 / 演出 / 全部类型 sheets (values only; no structural change).
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibitions)
$wsExhibit.Range("F2").Value = 195
$wsExhibit.Range("F4").Value = 405
$wsExhibit.Range("F5").Value = 973
$wsExhibit.Range("F6").Value = 5371
$wsExhibit.Range("F7").Value = 468
$wsExhibit.Range("F8").Value = 655
$wsExhibit.Range("F9").Value = 927
$wsExhibit.Range("F10").Value = 821
$wsExhibit.Range("F13").Value = 577
$wsExhibit.Range("F14").Value = 23
$wsExhibit.Range("F17").Value = 1791
$wsExhibit.Range("F18").Value = 1459
$wsExhibit.Range("F19").Value = 862
$wsExhibit.Range("F21").Value = 191
$wsExhibit.Range("F22").Value = 318
$wsExhibit.Range("F23").Value = 529
$wsExhibit.Range("F25").Value = 1049
$wsExhibit.Range("F27").Value = 524
$wsExhibit.Range("F28").Value = 2705
$wsExhibit.Range("F32").Value = 103
$wsExhibit.Range("F33").Value = 30
$wsExhibit.Range("F34").Value = 332
$wsExhibit.Range("F35").Value = 12
$wsExhibit.Range("F37").Value = 10
$wsExhibit.Range("F39").Value = 280
$wsExhibit.Range("F40").Value = 662
$wsExhibit.Range("F42").Value = 50
$wsExhibit.Range("F43").Value = 50
$wsExhibit.Range("F44").Value = 63

# 演出 (Performances)
$wsShow.Range("F4").Value = 162
$wsShow.Range("F6").Value = 111

# 全部类型 (All types, aggregated view)
$wsAll.Range("F3").Value = 195
$wsAll.Range("F5").Value = 973
$wsAll.Range("F7").Value = 5371
$wsAll.Range("F8").Value = 468
$wsAll.Range("F9").Value = 655
$wsAll.Range("F11").Value = 162
$wsAll.Range("F12").Value = 927
$wsAll.Range("F13").Value = 821
$wsAll.Range("F15").Value = 111
$wsAll.Range("F18").Value = 577
$wsAll.Range("F19").Value = 23
$wsAll.Range("F23").Value = 1791
$wsAll.Range("F24").Value = 1459
$wsAll.Range("F25").Value = 862
$wsAll.Range("F26").Value = 191
$wsAll.Range("F27").Value = 318
$wsAll.Range("F29").Value = 529
$wsAll.Range("F31").Value = 1049
$wsAll.Range("F32").Value = 2705
$wsAll.Range("F36").Value = 103
$wsAll.Range("F37").Value = 30
$wsAll.Range("F38").Value = 332
$wsAll.Range("F39").Value = 12
$wsAll.Range("F41").Value = 10
$wsAll.Range("F42").Value = 280
$wsAll.Range("F43").Value = 662
$wsAll.Range("F45").Value = 50
$wsAll.Range("F46").Value = 63
